$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = "Phase 01"
$ws.Range("B9").Value = "Department 01"
$ws.Range("C9").Value = 101
$ws.Range("D9").Value = "R1-1"
$ws.Range("E9").Value = "OFOI"
$ws.Range("F9").Value = "CAB00032"
$ws.Range("G9").Value = "Acid-Corrosive Cabinet, Manual Closing, Capacity: 30 ga"
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = "CRA-30"
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = "Eagle"
$ws.Range("L9").Value = 1
